$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the second row header cells with the 9 new labels (a..h, end)
$values = @("a", "b", "c", "d", "e", "f", "g", "h", "end")
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $values[$i]
}

# Move the active selection from E8 to G8
$ws.Range("G8").Select()
